# Update "epexspot_prices.xlsx"-style workbook with the latest day of data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column (AD) with hourly prices.
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting (font/border/alignment) of the previous header cell
# onto the new header cell before writing into it.
$wsPrix.Range("AC1").Copy()
$wsPrix.Range("AD1").PasteSpecial(-4122)
$wsPrix.Range("AD1").Value = "13-jul"

$prixValues = @{
    2  = 95.66
    3  = 87.25
    4  = 74.89
    5  = 60.72
    6  = 55.25
    7  = 54.07
    8  = 51.09
    9  = 62.31
    10 = 55
    11 = 22.5
    12 = 4.28
    13 = 0.05
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 9.37
    20 = 30.39
    21 = 70.1
    22 = 107.69
    23 = 122.8
    24 = 126.94
    25 = 114.4
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Range("AD$row").Value = $prixValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append a new row (27) for the next day.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force text so the "2025-07-11" date-looking string is not
# auto-converted into a date serial number, then drop the temporary
# number-format override so the cell keeps the sheet's default style.
$wsGaz.Range("A27").NumberFormat = "@"
$wsGaz.Range("A27").Value = "2025-07-11"
$wsGaz.Range("A27").Style = "Normal"
$wsGaz.Range("B27").Value = 34.8

# ---------------------------------------------------------------------
# Sheet "CO2": append a new row (27) for the next day.
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A27").NumberFormat = "@"
$wsCo2.Range("A27").Value = "2025-07-11"
$wsCo2.Range("A27").Style = "Normal"
$wsCo2.Range("B27").Value = 69.8
